$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.598.08'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.426.28'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.53'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.59'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +5.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +5.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.422.91'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +0.43%  '

$ws.Range("E10").Value = '  +2.56%  '

$ws.Range("E11").Value = '  +2.59%  '

$ws.Range("E12").Value = '  +1.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.021.90'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("E14").Value = '  +1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.24'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -1.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.529.69'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.93%  '

$ws.Range("E17").Value = '  +2.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.433.03'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("E19").Value = '  +0.95%  '

$ws.Range("E20").Value = '  +1.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '369.44'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +1.21%  '

$ws.Range("E22").Value = '  -0.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.24'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000126'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +6.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.534'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +2.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.86'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +2.99%  '

$ws.Range("E28").Value = '  +2.27%  '

$ws.Range("E29").Value = '  +0.33%  '

$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("E31").Value = '  +1.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.38'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -0.72%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("E34").Value = '  +1.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -0.50%  '

$ws.Range("E36").Value = '  +1.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.80'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +1.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.869'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -0.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.56'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -4.41%  '

$ws.Range("E40").Value = '  +3.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.66'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +3.37%  '

$ws.Range("E42").Value = '  +1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.711.98'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +0.51%  '

$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("E45").Value = '  +1.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.02'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +4.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.91'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '335.90'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +10.33%  '

$ws.Range("E49").Value = '  -0.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '32.38'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +7.51%  '

$ws.Range("E51").Value = '  +3.80%  '
